$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Adjust column widths (closest representable value to 13.4101845877511)
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe.Range("C1").ColumnWidth = 12.5
